$wb = $excel.ActiveWorkbook

# --- Rename the original sheet, add the new "31-05" sheet after it ---
$schedule = $wb.Worksheets.Item(1)
$schedule.Name = "Schedule"

$new = $wb.Worksheets.Add($null, $schedule)
$new.Name = "31-05"

# --- Populate the new sheet's data in the same order the shared-string
#     table in the target workbook was built (title last for row 4, the
#     row-6 item first, etc.) so unique-string ordering matches exactly ---
$new.Range("D6").Value = "Facility Bugs Insert Data"
$new.Range("E6").Value = "Hanya Item Ke 2 saya yang masuk Note Fasilitasnya"
$new.Range("D5").Value = "Name"
$new.Range("E5").Value = "Desc"
$new.Range("G5").Value = "Status"
$new.Range("F5").Value = "Root Cause/Analisa Sementara"
$new.Range("D4").Value = "INSERT NEW PACKAGE TOUR"
$new.Range("D7").Value = "Binding Image Caption"
$new.Range("D8").Value = "Add Time Travel"
$new.Range("E8").Value = "tambah waktu acara"
$new.Range("F8").Value = "dynamic datepicker"
$new.Range("E7").Value = "binding note image input ke image view"
$new.Range("D9").Value = "Set Default Video"
$new.Range("E9").Value = "setting default video ketika kosong"

# Numbering column (C6:C17 = 1..12)
$new.Range("C6").Value = 1
$new.Range("C7").Value = 2
$new.Range("C8").Value = 3
$new.Range("C9").Value = 4
$new.Range("C10").Value = 5
$new.Range("C11").Value = 6
$new.Range("C12").Value = 7
$new.Range("C13").Value = 8
$new.Range("C14").Value = 9
$new.Range("C15").Value = 10
$new.Range("C16").Value = 11
$new.Range("C17").Value = 12

# --- Title row formatting: merge D4:G4 and center it (new style) ---
$new.Range("D4:G4").Merge()
$new.Range("D4:G4").HorizontalAlignment = -4108

# --- Column widths ---
$new.Columns("C").ColumnWidth = 3.90625
$new.Columns("D").ColumnWidth = 20.6328125
$new.Columns("E").ColumnWidth = 43.90625
$new.Columns("F").ColumnWidth = 43.90625

# --- Selection / active tab on the new sheet ---
$new.Range("D10").Select()
$new.Activate()
